$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the entry for 11/17/2021 (row 84): course, hours, and notes
$ws.Range("B84").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C84").Value = 0.5
$ws.Range("D84").Value = "Finish 1 small problem"

# Move the active selection to C84, matching the saved view state
$ws.Range("C84").Select()
